$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# (losing the original formatting / introducing float rounding).
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '34.805.85'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = '1.828.95'
$ws.Range('E3').Value = '  +0.77%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = '230.72'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').Value = '39.48'
$ws.Range('E8').Value = '  -2.07%  '
$ws.Range('E9').Value = '  +1.20%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').Value = '0.0988'
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('D12').Value = '2.096.05'
$ws.Range('E12').Value = '  +0.96%  '
$ws.Range('D13').Value = '11.30'
$ws.Range('E13').Value = '  +1.63%  '
$ws.Range('D14').Value = '1.828.61'
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').Value = '0.666'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '4.64'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').Value = '34.813.65'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '69.38'
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').Value = '239.22'
$ws.Range('E20').Value = '  +0.40%  '
$ws.Range('D21').Value = '12.18'
$ws.Range('E21').Value = '  +2.87%  '
$ws.Range('D22').Value = '4.66'
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = '2.26'
$ws.Range('E24').Value = '  -0.57%  '
$ws.Range('D25').Value = '172.26'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').Value = '7.73'
$ws.Range('E26').Value = '  -1.39%  '
$ws.Range('D27').Value = '0.123'
$ws.Range('E27').Value = '  +2.57%  '
$ws.Range('D28').Value = '17.31'
$ws.Range('E28').Value = '  -0.75%  '
$ws.Range('D29').Value = '1.50'
$ws.Range('E29').Value = '  -8.13%  '
$ws.Range('E30').Value = '  +0.29%  '
$ws.Range('D31').Value = '0.0550'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').Value = '3.89'
$ws.Range('E32').Value = '  -0.73%  '
$ws.Range('D33').Value = '3.92'
$ws.Range('E33').Value = '  -1.38%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').Value = '1.84'
$ws.Range('E34').Value = '  +3.67%  '
$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D35').Value = '1.22'
$ws.Range('E35').Value = '  +7.40%  '
$ws.Range('E36').Value = '  +11.28%  '
$ws.Range('E37').Value = '  +2.53%  '
$ws.Range('D38').Value = '91.45'
$ws.Range('E38').Value = '  -1.78%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '1.04'
$ws.Range('E39').Value = '  +5.49%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.341.76'
$ws.Range('E40').Value = '  +2.37%  '
$ws.Range('D41').Value = '0.0193'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = '14.43'
$ws.Range('E42').Value = '  -1.89%  '
$ws.Range('E43').Value = '  -1.18%  '
$ws.Range('D44').Value = '2.24'
$ws.Range('E44').Value = '  -3.97%  '
$ws.Range('D45').Value = '2.75'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').Value = '6.27'
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('E47').Value = '  +1.85%  '
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('E50').Value = '  +4.11%  '
$ws.Range('D51').Value = '3.19'
$ws.Range('E51').Value = '  +12.77%  '

# Remove the temporary Text number-format again so the cell style matches
# the original (unstyled) cells while keeping the stored value as text.
$ws.Range('D5').ClearFormats()
$ws.Range('D8').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D15').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D18').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D22').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D26').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D28').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D32').ClearFormats()
$ws.Range('D33').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D35').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D39').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D42').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('D51').ClearFormats()
